$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 453. This shifts the existing rows
# 453-539 down to 454-540 (carrying their values/formatting with them,
# including the date-formatted style on column D), and grows the used
# range / dimension to A1:R540 automatically.
$ws.Rows.Item(453).EntireRow.Insert()

# Populate the newly-inserted row 453 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant/identical to the rest of
# this data block, so they are written as literals too (reading back
# `.Value` from a Range is not reliable in this environment, so we
# avoid round-tripping through it).
$ws.Range("A453").Value = 3
$ws.Range("B453").Value = "Femacal de La Calera"
$ws.Range("C453").Value = "Coquimbo"
$ws.Range("D453").Value = 45015
$ws.Range("E453").Value = 5
$ws.Range("F453").Value = 100112040
$ws.Range("G453").Value = "Cilantro"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 250
$ws.Range("K453").Value = 4000
$ws.Range("L453").Value = 4500
$ws.Range("M453").Value = 4240
$ws.Range("N453").Value = '$/docena de atados (3 kilos)'
$ws.Range("O453").Value = "Provincia de Quillota"
$ws.Range("P453").Value = 1413
$ws.Range("Q453").Value = 3
$ws.Range("R453").Value = "Hortaliza"
